$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the shrinking staircase matrix B2:K20.
# Each row's values shift one column to the left (dropping the old first
# value); rows 2-10 gain one brand-new trailing value, rows 11-20 simply
# lose their last (already-empty-after-shift) column.
$rowData = @{
    2  = @(1.56314828978338, -1.87688157821662, -1.93805557821662, 0.10340042178338, -0.47919457821662, -0.38351943421662, -0.64393857821662, 0.19151642178338, 0.09144346878338, -0.26234957821662)
    3  = @(-1.985496110432761, -2.046670110432761, -0.005214110432760821, -0.5878091104327609, -0.4921339664327609, -0.7525531104327607, 0.08290188956723921, -0.01717106343276081, -0.3709641104327608, -0.2227861104327608)
    4  = @(-1.838568389465022, 0.2028876105349782, -0.3797073894650218, -0.2840322454650218, -0.5444513894650218, 0.2910036105349783, 0.1909306575349782, -0.1628623894650218, -0.01468438946502179, 0.1427616105349782)
    5  = @(0.3193857809087786, -0.2632092190912214, -0.1675340750912214, -0.4279532190912214, 0.4075017809087786, 0.3074288279087786, -0.04636421909122137, 0.1018137809087786, 0.2592597809087786, -0.2414032190912214)
    6  = @(-0.2177158050940863, -0.1220406610940863, -0.3824598050940862, 0.4529951949059138, 0.3529222419059138, -0.0008708050940862688, 0.1473071949059137, 0.3047531949059137, -0.1959098050940863, 0.08445119490591374)
    7  = @(-0.1395947824540507, -0.4000139264540507, 0.4354410735459493, 0.3353681205459493, -0.01842492645405072, 0.1297530735459493, 0.2871990735459493, -0.2134639264540507, 0.06689707354594927, -0.1050779264540507)
    8  = @(-0.3119066075030444, 0.5235483924969555, 0.4234754394969555, 0.06968239249695551, 0.2178603924969555, 0.3753063924969555, -0.1253566075030445, 0.1550043924969555, -0.01697060750304449, 0.2888923924969555)
    9  = @(0.7021230259847391, 0.6020500729847391, 0.2482570259847391, 0.3964350259847391, 0.5538810259847391, 0.0532180259847391, 0.3335790259847391, 0.1616040259847391, 0.4674670259847391, -0.1448629740152609)
    10 = @(1.514077300737389, 1.160284253737389, 1.308462253737389, 1.465908253737389, 0.9652452537373891, 1.245606253737389, 1.073631253737389, 1.379494253737389, 0.7671642537373891, 1.455535253737389)
    11 = @(0.2163100177716323, 0.3644880177716323, 0.5219340177716323, 0.02127101777163229, 0.3016320177716323, 0.1296570177716323, 0.4355200177716323, -0.1768099822283677, 0.5115610177716323)
    12 = @(0.3684559152847414, 0.5259019152847414, 0.0252389152847414, 0.3055999152847414, 0.1336249152847414, 0.4394879152847414, -0.1728420847152586, 0.5155289152847414)
    13 = @(0.6615420054549828, 0.1608790054549828, 0.4412400054549828, 0.2692650054549828, 0.5751280054549828, -0.03720199454501721, 0.6511690054549828)
    14 = @(-0.07992717640068192, 0.2004338235993181, 0.02845882359931809, 0.3343218235993181, -0.2780081764006819, 0.4103628235993181)
    15 = @(0.1551026083472217, -0.01687239165277828, 0.2889906083472217, -0.3233393916527783, 0.3650316083472218)
    16 = @(-0.08373381358883215, 0.2221291864111679, -0.3902008135888321, 0.2981701864111679)
    17 = @(0.1925429840544512, -0.4197870159455487, 0.2685839840544513)
    18 = @(-0.4379376726751896, 0.2504333273248104)
    19 = @(0.2324015030761111)
    20 = @()
}

foreach ($r in 2..20) {
    $values = $rowData[$r]
    # Clear the full B:K span for this row first so columns beyond the
    # new shorter row length become blank again.
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 11)).ClearContents()
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $values[$i]
    }
}
